# Checkpoint before assistant change: Improve visibility of results by adjusting contrast for better readability
# Adds a new data row (row 2) to the "Visual Field", "CSV-1000", and "Pelli-Robinson" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Visual Field
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Visual Field")

$ws.Range("A2").Value = "ANI"

# Keep the test-date cell as literal text ("2025-05-27") instead of
# letting Excel auto-convert it into a date serial value.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-05-27"

$ws.Range("C2").Value = "2025-05-27T09:31:45.432Z"
$ws.Range("D2").Value = "2025-05-27T09:32:59.486Z"
$ws.Range("E2").Value = 74
$ws.Range("F2").Value = 48
$ws.Range("G2").Value = 22
$ws.Range("H2").Value = 45.83
$ws.Range("J2").Value = 48
$ws.Range("K2").Value = '[{''pointIndex'': 0, ''row'': 0, ''col'': 5, ''detected'': True, ''reactionTime'': 2632, ''timestamp'': ''2025-05-27T09:31:48.064Z''}, {''pointIndex'': 1, ''row'': 2, ''col'': 5, ''detected'': True, ''reactionTime'': 1368, ''timestamp'': ''2025-05-27T09:31:49.944Z''}, {''pointIndex'': 2, ''row'': 2, ''col'': 0, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:31:53.453Z''}, {''pointIndex'': 3, ''row'': 1, ''col'': 6, ''detected'': True, ''reactionTime'': 1784, ''timestamp'': ''2025-05-27T09:31:55.746Z''}, {''pointIndex'': 4, ''row'': 5, ''col'': 4, ''detected'': True, ''reactionTime'': 1856, ''timestamp'': ''2025-05-27T09:31:58.118Z''}, {''pointIndex'': 5, ''row'': 4, ''col'': 3, ''detected'': True, ''reactionTime'': 1200, ''timestamp'': ''2025-05-27T09:31:59.831Z''}, {''pointIndex'': 6, ''row'': 0, ''col'': 2, ''detected'': True, ''reactionTime'': 1316, ''timestamp'': ''2025-05-27T09:32:01.649Z''}, {''pointIndex'': 7, ''row'': 1, ''col'': 1, ''detected'': True, ''reactionTime'': 1185, ''timestamp'': ''2025-05-27T09:32:03.344Z''}, {''pointIndex'': 8, ''row'': 3, ''col'': 0, ''detected'': True, ''reactionTime'': 2535, ''timestamp'': ''2025-05-27T09:32:06.383Z''}, {''pointIndex'': 9, ''row'': 4, ''col'': 6, ''detected'': True, ''reactionTime'': 1928, ''timestamp'': ''2025-05-27T09:32:08.816Z''}, {''pointIndex'': 10, ''row'': 1, ''col'': 3, ''detected'': True, ''reactionTime'': 1656, ''timestamp'': ''2025-05-27T09:32:10.983Z''}, {''pointIndex'': 11, ''row'': 5, ''col'': 2, ''detected'': True, ''reactionTime'': 1385, ''timestamp'': ''2025-05-27T09:32:12.883Z''}, {''pointIndex'': 12, ''row'': 1, ''col'': 2, ''detected'': True, ''reactionTime'': 1139, ''timestamp'': ''2025-05-27T09:32:14.533Z''}, {''pointIndex'': 13, ''row'': 2, ''col'': 7, ''detected'': True, ''reactionTime'': 1669, ''timestamp'': ''2025-05-27T09:32:16.714Z''}, {''pointIndex'': 14, ''row'': 1, ''col'': 5, ''detected'': True, ''reactionTime'': 1116, ''timestamp'': ''2025-05-27T09:32:18.332Z''}, {''pointIndex'': 15, ''row'': 0, ''col'': 5, ''detected'': True, ''reactionTime'': 1957, ''timestamp'': ''2025-05-27T09:32:20.799Z''}, {''pointIndex'': 16, ''row'': 5, ''col'': 0, ''detected'': True, ''reactionTime'': 1701, ''timestamp'': ''2025-05-27T09:32:23.012Z''}, {''pointIndex'': 17, ''row'': 0, ''col'': 1, ''detected'': True, ''reactionTime'': 1619, ''timestamp'': ''2025-05-27T09:32:25.131Z''}, {''pointIndex'': 18, ''row'': 4, ''col'': 7, ''detected'': True, ''reactionTime'': 1590, ''timestamp'': ''2025-05-27T09:32:27.235Z''}, {''pointIndex'': 19, ''row'': 3, ''col'': 5, ''detected'': True, ''reactionTime'': 1267, ''timestamp'': ''2025-05-27T09:32:29.018Z''}, {''pointIndex'': 20, ''row'': 3, ''col'': 0, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:30.385Z''}, {''pointIndex'': 21, ''row'': 0, ''col'': 2, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:30.940Z''}, {''pointIndex'': 22, ''row'': 5, ''col'': 7, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:31.488Z''}, {''pointIndex'': 23, ''row'': 3, ''col'': 2, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:32.101Z''}, {''pointIndex'': 24, ''row'': 2, ''col'': 2, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:32.998Z''}, {''pointIndex'': 25, ''row'': 2, ''col'': 7, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:33.640Z''}, {''pointIndex'': 26, ''row'': 0, ''col'': 1, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:34.353Z''}, {''pointIndex'': 27, ''row'': 1, ''col'': 6, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:37.882Z''}, {''pointIndex'': 28, ''row'': 5, ''col'': 3, ''detected'': True, ''reactionTime'': 1304, ''timestamp'': ''2025-05-27T09:32:39.701Z''}, {''pointIndex'': 29, ''row'': 0, ''col'': 3, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:40.707Z''}, {''pointIndex'': 30, ''row'': 2, ''col'': 0, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:41.305Z''}, {''pointIndex'': 31, ''row'': 4, ''col'': 3, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:41.882Z''}, {''pointIndex'': 32, ''row'': 5, ''col'': 0, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:42.745Z''}, {''pointIndex'': 33, ''row'': 2, ''col'': 5, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:43.322Z''}, {''pointIndex'': 34, ''row'': 1, ''col'': 0, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:43.907Z''}, {''pointIndex'': 35, ''row'': 1, ''col'': 5, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:44.448Z''}, {''pointIndex'': 36, ''row'': 4, ''col'': 5, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:45.090Z''}, {''pointIndex'': 37, ''row'': 3, ''col'': 3, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:45.746Z''}, {''pointIndex'': 38, ''row'': 5, ''col'': 1, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:46.407Z''}, {''pointIndex'': 39, ''row'': 1, ''col'': 5, ''detected'': True, ''reactionTime'': 1789, ''timestamp'': ''2025-05-27T09:32:48.703Z''}, {''pointIndex'': 40, ''row'': 5, ''col'': 7, ''detected'': True, ''reactionTime'': 1444, ''timestamp'': ''2025-05-27T09:32:50.653Z''}, {''pointIndex'': 41, ''row'': 0, ''col'': 7, ''detected'': False, ''reactionTime'': 921, ''timestamp'': ''2025-05-27T09:32:52.081Z''}, {''pointIndex'': 42, ''row'': 1, ''col'': 5, ''detected'': False, ''reactionTime'': 857, ''timestamp'': ''2025-05-27T09:32:53.439Z''}, {''pointIndex'': 43, ''row'': 3, ''col'': 1, ''detected'': False, ''reactionTime'': 524, ''timestamp'': ''2025-05-27T09:32:54.470Z''}, {''pointIndex'': 44, ''row'': 0, ''col'': 0, ''detected'': False, ''reactionTime'': 701, ''timestamp'': ''2025-05-27T09:32:55.684Z''}, {''pointIndex'': 45, ''row'': 2, ''col'': 6, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:57.621Z''}, {''pointIndex'': 46, ''row'': 3, ''col'': 5, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:58.299Z''}, {''pointIndex'': 47, ''row'': 5, ''col'': 0, ''detected'': False, ''reactionTime'': 3000, ''timestamp'': ''2025-05-27T09:32:58.983Z''}]'
$ws.Range("L2").Value = 22

# ---------------------------------------------------------------
# Sheet: CSV-1000
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CSV-1000")

$ws.Range("A2").Value = "ANI"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-05-27"

$ws.Range("C2").Value = "2025-05-27T09:33:55.020Z"
$ws.Range("D2").Value = "2025-05-27T09:34:38.943Z"
$ws.Range("E2").Value = 44
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 60
$ws.Range("J2").Value = "English"
$ws.Range("K2").Value = '[100, 50, 25, 12.5, 6.25, 3.125, 1.56, 0.78, 0.39, 0.195]'
$ws.Range("L2").Value = 12

# ---------------------------------------------------------------
# Sheet: Pelli-Robinson
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pelli-Robinson")

$ws.Range("A2").Value = "ANI"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-05-27"

$ws.Range("C2").Value = "2025-05-27T09:35:34.132Z"
$ws.Range("D2").Value = "2025-05-27T09:38:55.761Z"
$ws.Range("E2").Value = 202
$ws.Range("F2").Value = 40
$ws.Range("G2").Value = 16
$ws.Range("H2").Value = 40
